# Redesigned phase envelope utilities interface
# Update the electrolysis power input on the Input sheet; all dependent
# formulas on Output and Calculations recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
$ws.Range("B12").Value = 15.419249390939662

$excel.CalculateFullRebuild()
